$wb = $excel.ActiveWorkbook

# --- Sheets 1-4: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
#     "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)"
#     All share the same A-column layout (Fonte/Tecnologia header + rows 2-12).
$ws1 = $wb.Worksheets.Item(1)
$fonteSheets = @(1, 2, 3, 4)
foreach ($idx in $fonteSheets) {
    $ws = $wb.Worksheets.Item($idx)

    # Add header cell A1 "Fonte/Tecnologia", matching the style of B1 (bold/centered/bordered)
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # Fix accents in rows 2-12
    $ws.Range("A2").Value = "Hidro"
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A5").Value = "Nuclear"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A7").Value = "Biomassa"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A9").Value = "Solar"
    $ws.Range("A10").Value = "Outros"
    $ws.Range("A11").Value = "Pot. Compl."
    $ws.Range("A12").Value = "GD"

    # Remove the header style from A2:A12 (now plain text cells, no style)
    $ws.Range("A2:A12").ClearFormats()
}

# --- Sheet 5: "Emissoes Totais (MtCO2eq)"
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A2:A3").ClearFormats()

# Remove row 4 ("Teto") entirely
$ws5.Rows(4).Delete()

# --- Sheet 6: "Custo Total (bilhões de R$)"
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$ws6.Range("A1").Value = "Tipo Expansão"

# B1 header text changes from "Custo" to "2015" (kept as text, not a number)
$ws6.Range("B1").NumberFormat = "@"
$ws6.Range("B1").Value = "2015"
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4122)

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 588
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
$ws6.Range("A2:A3").ClearFormats()

$excel.CutCopyMode = 0
